$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds Property/Value pairs; row 7 is the "Experimental" row.
$ws = $wb.Worksheets.Item("Metadata")

# ValueSet now carries the (previously-missing) required "experimental"
# boolean element -- record it as text "true" in the Value column.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# The export's "Date" metadata value also moved forward to reflect the
# re-run that added the experimental flag.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
